{"js": "// 1) \"Inspect the Sample-Live-Sync directory hierarchy ...\" paragraph:\n//    - \"Sample-Live-Sync\" -> \"CloudSdkSyncSample\"\n//    - insert \" (except Sample-Live-Sync.*\" right after \"pdb\"\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfunction findParagraphByText(items, needle) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(needle) !== -1) {\n      return items[i];\n    }\n  }\n  return null;\n}\n\nconst inspectPara = findParagraphByText(\n  paragraphs.items,\n  \"Inspect the Sample-Live-Sync directory hierarchy\"\n);\nif (!inspectPara) {\n  throw new Error(\"Could not find the 'Inspect the Sample-Live-Sync...' paragraph\");\n}\n\nconst sampleLiveSyncHits = inspectPara.search(\"Sample-Live-Sync\", { matchCase: true });\nsampleLiveSyncHits.load(\"text\");\nawait context.sync();\nsampleLiveSyncHits.items[0].insertText(\"CloudSdkSyncSample\", \"Replace\");\nawait context.sync();\n\nconst pdbHits = inspectPara.search(\"pdb\", { matchCase: true });\npdbHits.load(\"text\");\nawait context.sync();\npdbHits.items[0].insertText(\"pdb (except Sample-Live-Sync.*\", \"Replace\");\nawait context.sync();\n\n// 2) Branch/version string \"20130218A0Release0_1_6\" -> \"20130218A0Release0.1.6\"\n//    (occurs 3 times: git branch, git checkout, git push -u origin)\nconst versionHits = context.document.body.search(\"20130218A0Release0_1_6\", {\n  matchCase: true\n});\nversionHits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < versionHits.items.length; i++) {\n  versionHits.items[i].insertText(\"20130218A0Release0.1.6\", \"Replace\");\n}\nawait context.sync();\n\n// 3) Move the \"_GoBack\" bookmark from the \"Build Sample-Live-Sync project...\"\n//    paragraph to the start of the \"Git branch - -d <branch>\" paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\nconst deleteBranchPara = findParagraphByText(\n  paragraphs2.items,\n  \"Git branch - -d <branch>\"\n);\nif (!deleteBranchPara) {\n  throw new Error(\"Could not find the 'Git branch - -d <branch>' paragraph\");\n}\nconst startOfDeleteBranchPara = deleteBranchPara.getRange(\"Start\");\nstartOfDeleteBranchPara.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word COM interop script - apply the \"staging release\" document edits.\n$d = $word.ActiveDocument\n\n# --- helper: find the first paragraph whose text contains $needle ---\nfunction Find-ParagraphRange($doc, $needle) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text -like \"*$needle*\") {\n            return $p.Range\n        }\n    }\n    return $null\n}\n\n# 1) \"Inspect the Sample-Live-Sync directory hierarchy ...\" paragraph:\n#    - \"Sample-Live-Sync\" -> \"CloudSdkSyncSample\"\n#    - insert \" (except Sample-Live-Sync.*\" right after \"pdb\"\n$inspectRange = Find-ParagraphRange $d \"Inspect the Sample-Live-Sync directory hierarchy\"\n$find1 = $inspectRange.Find\n$find1.Execute(\"Sample-Live-Sync\", $false, $false, $false, $false, $false, $true, 1, $false, \"CloudSdkSyncSample\", 2)\n\n$inspectRange2 = Find-ParagraphRange $d \"CloudSdkSyncSample directory hierarchy\"\n$find2 = $inspectRange2.Find\n$find2.Execute(\"pdb\", $false, $false, $false, $false, $false, $true, 1, $false, \"pdb (except Sample-Live-Sync.*\", 2)\n\n# 2) Branch/version string \"20130218A0Release0_1_6\" -> \"20130218A0Release0.1.6\"\n#    (occurs 3 times: git branch, git checkout, git push -u origin)\n$versionFind = $d.Content.Find\n$versionFind.Execute(\"20130218A0Release0_1_6\", $false, $false, $false, $false, $false, $true, 1, $false, \"20130218A0Release0.1.6\", 2)\n\n# 3) Move the \"_GoBack\" bookmark from the \"Build Sample-Live-Sync project...\"\n#    paragraph to the start of the \"Git branch - -d <branch>\" paragraph.\n$deleteBranchRange = Find-ParagraphRange $d \"Git branch - -d <branch>\"\n$startRange = $d.Range($deleteBranchRange.Start, $deleteBranchRange.Start)\n$d.Bookmarks.Add(\"_GoBack\", $startRange)\n"}
